$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first six data rows (old rows 2-7); everything below shifts up.
$ws.Rows("2:7").Delete()

# New accelerometer samples appended at the bottom (the window slid forward).
$newRows = @(
    @(-12.76972269274368, -22.28177534673628, 15.843390720408),
    @(-17.05221034571044, -9.689534740349583, -7.262813420639774),
    @(7.601734495654355, -7.369463099647012, 0.05982191538072679),
    @(-5.084137159524504, -17.42479633056011, 13.42315847357525),
    @(-2.988111949458581, -1.541268535496201, 24.50189582588753),
    @(16.88818173740272, -21.43393101642941, 26.25781544950798),
    @(11.71483505878375, -29.24488582807863, 14.33557478914648),
    @(-17.89474326064904, 1.058929585918881, -6.072563363104724),
    @(-6.588068613071894, -42.04801777711909, 27.395873059932),
    @(4.73087814911125, -46.72824330182412, 40.59463504909236),
    @(13.50610577691452, -5.020019875359202, 23.69000314928804),
    @(-13.47701175925547, -14.91354519067346, 16.03211706692484),
    @(-16.10707013631821, -7.06407377891983, -12.22038123533914),
    @(17.60490359473454, -9.783274660405542, -4.873167308335432),
    @(-13.56971339589527, -42.06590333918953, 1.837090728209249),
    @(-0.4818755535735217, -32.3410521585916, 8.205087691238255)
)

$startRow = 16
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
